$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1
$ws.Range("H1").Value = "Save"

# Copy formatting (bold font, thin border, centered/top alignment) from the
# existing "sum" header (G1) onto the new "Save" header so it matches the
# other header cells exactly.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Fill in the "Save" column values for the data rows
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 1
